$wb = $excel.ActiveWorkbook

# Hunk 0: sheet ALC, row 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 552.8570999999999
$ws.Range("J2").Value = 733.5
$ws.Range("L2").Value = 733.5
$ws.Range("N2").Value = -959.5

# Hunk 1: sheet ALC, row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1850
$ws.Range("J17").Value = 1850
$ws.Range("L17").Value = 5550
$ws.Range("N17").Value = -5886

# Hunk 2: sheet ALC, row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3250
$ws.Range("I32").Value = 2500
$ws.Range("J32").Value = 4000
$ws.Range("K32").Value = 2500
$ws.Range("L32").Value = 4000
$ws.Range("M32").Value = -2174
$ws.Range("N32").Value = -4652

# Hunk 3: sheet ALC, row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1613.3636
$ws.Range("I40").Value = 1699.875
$ws.Range("J40").Value = 1382.6666
$ws.Range("K40").Value = 1699.875
$ws.Range("L40").Value = 1382.6666
$ws.Range("M40").Value = -1524.875
$ws.Range("N40").Value = -1732.6666

# Hunk 4: sheet ALC, row 53
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 356.27274
$ws.Range("I53").Value = 402.5
$ws.Range("J53").Value = 233
$ws.Range("K53").Value = 402.5
$ws.Range("L53").Value = 233
$ws.Range("M53").Value = 234.5
$ws.Range("N53").Value = -1507

# Hunk 5: sheet ALC, row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 14583.6
$ws.Range("I113").Value = 13900
$ws.Range("J113").Value = 14754.5
$ws.Range("K113").Value = 13900
$ws.Range("L113").Value = 14754.5
$ws.Range("M113").Value = -10646
$ws.Range("N113").Value = -21262.5

# Hunk 6: sheet ALC, row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

# Hunk 7: sheet ARM, row 46
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 19333
$ws.Range("I46").Value = 19333
$ws.Range("K46").Value = 19333
$ws.Range("M46").Value = -19014

# Hunk 8: sheet ARM, row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 335319.66
$ws.Range("J102").Value = 2975
$ws.Range("L102").Value = 2975
$ws.Range("N102").Value = -6219

# Hunk 9: sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 13333.333
$ws.Range("I122").Value = 13333.333
$ws.Range("K122").Value = 39999.999
$ws.Range("M122").Value = -37549.999

# Hunk 10: sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5500
$ws.Range("I132").Value = 5500
$ws.Range("K132").Value = 16500
$ws.Range("M132").Value = -13970

# Hunk 11: sheet BSM, row 45
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

# Hunk 12: sheet BSM, row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1600
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# Hunk 13: sheet BSM, row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1600
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# Hunk 14: sheet BSM, row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3031.25
$ws.Range("I94").Value = 3325
$ws.Range("K94").Value = 3325
$ws.Range("M94").Value = -2874

# Hunk 15: sheet BSM, row 130
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 139999.67
$ws.Range("J130").Value = 139999.67
$ws.Range("L130").Value = 139999.67
$ws.Range("N130").Value = -150039.67

# Hunk 16: sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5560.143
$ws.Range("J134").Value = 5130
$ws.Range("L134").Value = 15390
$ws.Range("N134").Value = -20460

# Hunk 17: sheet CRP, row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 57.866665
$ws.Range("I7").Value = 69.52631
$ws.Range("K7").Value = 69.52631
$ws.Range("M7").Value = 43.47369

# Hunk 18: sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2547.1428
$ws.Range("I31").Value = 2547.1428
$ws.Range("K31").Value = 2547.1428
$ws.Range("M31").Value = -2252.1428

# Hunk 19: sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2547.1428
$ws.Range("I34").Value = 2547.1428
$ws.Range("K34").Value = 2547.1428
$ws.Range("M34").Value = -2345.1428

# Hunk 20: sheet CRP, row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1027.8334
$ws.Range("J94").Value = 1035.4
$ws.Range("L94").Value = 1035.4
$ws.Range("N94").Value = -1937.4

# Hunk 21: sheet CRP, row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 502718.4
$ws.Range("I99").Value = 1624.75
$ws.Range("J99").Value = 836780.8
$ws.Range("K99").Value = 1624.75
$ws.Range("L99").Value = 836780.8
$ws.Range("M99").Value = -126.75
$ws.Range("N99").Value = -839776.8

# Hunk 22: sheet CRP, row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 502718.4
$ws.Range("I126").Value = 1624.75
$ws.Range("J126").Value = 836780.8
$ws.Range("K126").Value = 4874.25
$ws.Range("L126").Value = 2510342.4
$ws.Range("M126").Value = -2404.25
$ws.Range("N126").Value = -2515282.4

# Hunk 23: sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3977.4375
$ws.Range("I132").Value = 3088.7693
$ws.Range("J132").Value = 7828.3335
$ws.Range("K132").Value = 9266.3079
$ws.Range("L132").Value = 23485.0005
$ws.Range("M132").Value = -6736.3079
$ws.Range("N132").Value = -28545.0005

# Hunk 24: sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1210.5
$ws.Range("J134").Value = 1210
$ws.Range("L134").Value = 3630
$ws.Range("N134").Value = -8700

# Hunk 25: sheet CUL, row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 177
$ws.Range("J2").Value = 435.16666
$ws.Range("L2").Value = 2610.99996
$ws.Range("N2").Value = -2836.99996

# Hunk 26: sheet CUL, row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 62199.562
$ws.Range("I4").Value = 63459.13
$ws.Range("K4").Value = 190377.39
$ws.Range("M4").Value = -190265.39

# Hunk 27: sheet CUL, row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 191.66667
$ws.Range("I38").Value = 191.66667
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 575.00001
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -228.00001
$ws.Range("N38").ClearContents()

# Hunk 28: sheet CUL, row 49
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 3999.5
$ws.Range("J49").Value = 3999.5
$ws.Range("L49").Value = 11998.5
$ws.Range("N49").Value = -12310.5

# Hunk 29: sheet CUL, row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 51449
$ws.Range("I139").Value = 2898
$ws.Range("K139").Value = 8694
$ws.Range("M139").Value = -3554

# Hunk 30: sheet GSM, row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1698.8334
$ws.Range("J107").Value = 1698.8334
$ws.Range("L107").Value = 1698.8334
$ws.Range("N107").Value = -5538.8334

# Hunk 31: sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4501.4375
$ws.Range("I122").Value = 3560
$ws.Range("K122").Value = 10680
$ws.Range("M122").Value = -8230

# Hunk 32: sheet LTW, row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7000
$ws.Range("I7").Value = 7000
$ws.Range("K7").Value = 7000
$ws.Range("M7").Value = -6888

# Hunk 33: sheet LTW, row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1737.5
$ws.Range("I16").Value = 1077.7778
$ws.Range("J16").Value = 7675
$ws.Range("K16").Value = 1077.7778
$ws.Range("L16").Value = 7675
$ws.Range("M16").Value = -907.7778000000001
$ws.Range("N16").Value = -8015

# Hunk 34: sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3430.6924
$ws.Range("J122").Value = 3685.5715
$ws.Range("L122").Value = 11056.7145
$ws.Range("N122").Value = -15956.7145

# Hunk 35: sheet LTW, row 124
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H124").Value = 15214.5
$ws.Range("J124").Value = 15214.5
$ws.Range("L124").Value = 15214.5
$ws.Range("N124").Value = -25034.5

# Hunk 36: sheet LTW, row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 7000
$ws.Range("I126").Value = 7000
$ws.Range("K126").Value = 21000
$ws.Range("M126").Value = -18530

# Hunk 37: sheet WVR, row 10
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 2755.5
$ws.Range("I10").Value = 11
$ws.Range("J10").Value = 5500
$ws.Range("K10").Value = 11
$ws.Range("L10").Value = 5500
$ws.Range("M10").Value = 158
$ws.Range("N10").Value = -5838

# Hunk 38: sheet WVR, row 14
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1114222.2
$ws.Range("J14").Value = 3500
$ws.Range("L14").Value = 3500
$ws.Range("N14").Value = -3836

# Hunk 39: sheet WVR, row 80
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 15000
$ws.Range("J80").Value = 15000
$ws.Range("L80").Value = 15000
$ws.Range("N80").Value = -16996

# Hunk 40: sheet WVR, row 83
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H83").Value = 15000
$ws.Range("J83").Value = 15000
$ws.Range("L83").Value = 45000
$ws.Range("N83").Value = -54984

# Hunk 41: sheet WVR, row 101
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 23068
$ws.Range("J101").Value = 23068
$ws.Range("L101").Value = 23068
$ws.Range("N101").Value = -29558

# Hunk 42: sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 414.5
$ws.Range("I122").Value = 414.5
$ws.Range("K122").Value = 1243.5
$ws.Range("M122").Value = 1206.5

# Hunk 43: sheet WVR, row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3485
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 3485
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 10455
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -15395

# Hunk 44: sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1314.5555
$ws.Range("I132").Value = 1023.6667
$ws.Range("K132").Value = 3071.0001
$ws.Range("M132").Value = -541.0001000000002
